# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 176.73685
$ws.Range("I33").Value = 186.64706
$ws.Range("J33").Value = 92.5
$ws.Range("K33").Value = 186.64706
$ws.Range("L33").Value = 92.5
$ws.Range("M33").Value = 42.35293999999999
$ws.Range("N33").Value = -550.5
$ws.Range("H64").Value = 4123.5293
$ws.Range("J64").Value = 4240
$ws.Range("L64").Value = 4240
$ws.Range("N64").Value = -4736
$ws.Range("H67").Value = 4123.5293
$ws.Range("J67").Value = 4240
$ws.Range("L67").Value = 4240
$ws.Range("N67").Value = -5956
$ws.Range("H98").Value = 3089534.8
$ws.Range("I98").Value = 3970827
$ws.Range("J98").Value = 5011.875
$ws.Range("K98").Value = 3970827
$ws.Range("L98").Value = 5011.875
$ws.Range("M98").Value = -3969329
$ws.Range("N98").Value = -8007.875
$ws.Range("H111").Value = 2586.5715
$ws.Range("I111").Value = 2765.8
$ws.Range("K111").Value = 8297.400000000001
$ws.Range("M111").Value = -5230.400000000001
$ws.Range("H122").Value = 3089534.8
$ws.Range("I122").Value = 3970827
$ws.Range("J122").Value = 5011.875
$ws.Range("K122").Value = 11912481
$ws.Range("L122").Value = 15035.625
$ws.Range("M122").Value = -11910031
$ws.Range("N122").Value = -19935.625
$ws.Range("H125").Value = 805.41174
$ws.Range("I125").Value = 707
$ws.Range("J125").Value = 859.0909
$ws.Range("K125").Value = 6363
$ws.Range("L125").Value = 7731.8181
$ws.Range("M125").Value = -3903
$ws.Range("N125").Value = -12651.8181
$ws.Range("H137").Value = 1067.6489
$ws.Range("I137").Value = 892.0968
$ws.Range("J137").Value = 1407.7812
$ws.Range("K137").Value = 2676.2904
$ws.Range("L137").Value = 4223.3436
$ws.Range("M137").Value = -126.2903999999999
$ws.Range("N137").Value = -9323.3436
$ws.Range("H138").Value = 1949.01
$ws.Range("I138").Value = 937
$ws.Range("J138").Value = 2681.8447
$ws.Range("K138").Value = 2811
$ws.Range("L138").Value = 8045.534100000001
$ws.Range("M138").Value = 2329
$ws.Range("N138").Value = -18325.5341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15465.805
$ws.Range("I32").Value = 13556.172
$ws.Range("J32").Value = 21541.908
$ws.Range("K32").Value = 13556.172
$ws.Range("L32").Value = 21541.908
$ws.Range("M32").Value = -13269.172
$ws.Range("N32").Value = -22115.908
$ws.Range("H74").Value = 909.3333
$ws.Range("I74").Value = 768
$ws.Range("J74").Value = 1333.3334
$ws.Range("K74").Value = 768
$ws.Range("L74").Value = 1333.3334
$ws.Range("M74").Value = 106
$ws.Range("N74").Value = -3081.3334
$ws.Range("H77").Value = 909.3333
$ws.Range("I77").Value = 768
$ws.Range("J77").Value = 1333.3334
$ws.Range("K77").Value = 3840
$ws.Range("L77").Value = 6666.666999999999
$ws.Range("M77").Value = 528
$ws.Range("N77").Value = -15402.667
$ws.Range("H88").Value = 37861.715
$ws.Range("J88").Value = 37861.715
$ws.Range("L88").Value = 37861.715
$ws.Range("N88").Value = -38673.715
$ws.Range("H91").Value = 37861.715
$ws.Range("J91").Value = 37861.715
$ws.Range("L91").Value = 37861.715
$ws.Range("N91").Value = -40669.715
$ws.Range("H122").Value = 40000710
$ws.Range("I122").Value = 41667380
$ws.Range("J122").Value = 666
$ws.Range("K122").Value = 125002140
$ws.Range("L122").Value = 1998
$ws.Range("M122").Value = -124999690
$ws.Range("N122").Value = -6898

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6956132.5
$ws.Range("I134").Value = 7764622.5
$ws.Range("J134").Value = 3120
$ws.Range("K134").Value = 23293867.5
$ws.Range("L134").Value = 9360
$ws.Range("M134").Value = -23291332.5
$ws.Range("N134").Value = -14430

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1145.05
$ws.Range("I16").Value = 1111.3125
$ws.Range("J16").Value = 1280
$ws.Range("K16").Value = 1111.3125
$ws.Range("L16").Value = 1280
$ws.Range("M16").Value = -824.3125
$ws.Range("N16").Value = -1854
$ws.Range("H94").Value = 2756.5
$ws.Range("I94").Value = 10000
$ws.Range("J94").Value = 2273.6
$ws.Range("K94").Value = 10000
$ws.Range("L94").Value = 2273.6
$ws.Range("M94").Value = -9549
$ws.Range("N94").Value = -3175.6
$ws.Range("H99").Value = 250001100
$ws.Range("I99").Value = 333334140
$ws.Range("K99").Value = 333334140
$ws.Range("M99").Value = -333332642
$ws.Range("H113").Value = 1145.05
$ws.Range("I113").Value = 1111.3125
$ws.Range("J113").Value = 1280
$ws.Range("K113").Value = 1111.3125
$ws.Range("L113").Value = 1280
$ws.Range("M113").Value = 1058.6875
$ws.Range("N113").Value = -5620
$ws.Range("H122").Value = 28572650
$ws.Range("I122").Value = 35715310
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 107145930
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -107143480
$ws.Range("N122").Value = -10897
$ws.Range("H126").Value = 250001100
$ws.Range("I126").Value = 333334140
$ws.Range("K126").Value = 1000002420
$ws.Range("M126").Value = -999999950
$ws.Range("H132").Value = 8551157
$ws.Range("I132").Value = 13334205
$ws.Range("J132").Value = 10000.429
$ws.Range("K132").Value = 40002615
$ws.Range("L132").Value = 30001.287
$ws.Range("M132").Value = -40000085
$ws.Range("N132").Value = -35061.287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1350.0588
$ws.Range("I5").Value = 415.8889
$ws.Range("J5").Value = 2401
$ws.Range("K5").Value = 1247.6667
$ws.Range("L5").Value = 7203
$ws.Range("M5").Value = -1135.6667
$ws.Range("N5").Value = -7427
$ws.Range("H125").Value = 3357.1428
$ws.Range("J125").Value = 4380
$ws.Range("L125").Value = 13140
$ws.Range("N125").Value = -22980
$ws.Range("H132").Value = 1735.3765
$ws.Range("I132").Value = 614.8
$ws.Range("J132").Value = 2080.1692
$ws.Range("K132").Value = 5533.2
$ws.Range("L132").Value = 18721.5228
$ws.Range("M132").Value = -3003.2
$ws.Range("N132").Value = -23781.5228
$ws.Range("H133").Value = 14081.818
$ws.Range("I133").Value = 8940
$ws.Range("K133").Value = 26820
$ws.Range("M133").Value = -21760
$ws.Range("H135").Value = 1350.0588
$ws.Range("I135").Value = 415.8889
$ws.Range("J135").Value = 2401
$ws.Range("K135").Value = 3743.0001
$ws.Range("L135").Value = 21609
$ws.Range("M135").Value = -1208.0001
$ws.Range("N135").Value = -26679

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 5902.4443
$ws.Range("J57").Value = 6015.25
$ws.Range("L57").Value = 6015.25
$ws.Range("N57").Value = -7655.25
$ws.Range("H122").Value = 45458004
$ws.Range("I122").Value = 83334920
$ws.Range("J122").Value = 5709.8
$ws.Range("K122").Value = 250004760
$ws.Range("L122").Value = 17129.4
$ws.Range("M122").Value = -250002310
$ws.Range("N122").Value = -22029.4
$ws.Range("H123").Value = 16123.556
$ws.Range("J123").Value = 16123.556
$ws.Range("L123").Value = 16123.556
$ws.Range("N123").Value = -21023.556
$ws.Range("H132").Value = 15169523
$ws.Range("I132").Value = 18888352
$ws.Range("J132").Value = 8139.4614
$ws.Range("K132").Value = 56665056
$ws.Range("L132").Value = 24418.3842
$ws.Range("M132").Value = -56662526
$ws.Range("N132").Value = -29478.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11365582
$ws.Range("I122").Value = 1873.4667
$ws.Range("J122").Value = 35716384
$ws.Range("K122").Value = 5620.4001
$ws.Range("L122").Value = 107149152
$ws.Range("M122").Value = -3170.4001
$ws.Range("N122").Value = -107154052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1054.6666
$ws.Range("I122").Value = 1054.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3163.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -713.9998000000001
$ws.Range("N122").ClearContents()
